$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Components added: update the two input values that drive the
# resistor-divider calculations (downstream formula cells recalc
# automatically).
$ws.Range("H2").Value = 6800
$ws.Range("B6").Value = 120000

# First schematic drawn / cell focus moved to H3.
$ws.Range("H3").Select()
